$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.AddChart2(-1, 74, 100, 100, 300, 300)
$chart = $shp.Chart
$s = $chart.SeriesCollection().NewSeries()
$s.XValues = $ws.Range("J159:J166")
$s.Values = $ws.Range("K159:K166")

$axCat = $chart.Axes(1, 1)
Write-Output "catAx ok"
$axVal = $chart.Axes(2, 1)
Write-Output "valAx ok"
